# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46075
$ws.Range("C2").Value = 27.23
$ws.Range("D2").Value = 22.7
$ws.Range("E2").Value = 15.27
$ws.Range("F2").Value = 9.390000000000001
$ws.Range("G2").Value = 9.18
$ws.Range("H2").Value = 11.74
$ws.Range("I2").Value = 19.58
$ws.Range("J2").Value = 6.23
$ws.Range("K2").Value = 1.22
$ws.Range("L2").Value = 0.32
$ws.Range("M2").Value = 0.23
$ws.Range("N2").Value = 0.54
$ws.Range("O2").Value = 0.22
$ws.Range("P2").Value = 0.19
$ws.Range("Q2").Value = 0.17
$ws.Range("R2").Value = 0.28
$ws.Range("S2").Value = 6.63
$ws.Range("T2").Value = 20.61
$ws.Range("U2").Value = 32.18
$ws.Range("V2").Value = 38.81
$ws.Range("W2").Value = 73
$ws.Range("X2").Value = 45.46
$ws.Range("Y2").Value = 33.04
$ws.Range("Z2").Value = 16.97
$ws.Range("AB2").Value = 47.58
$ws.Range("AD2").Value = 55.9
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 39.25
$ws.Range("AG2").Value = "3h-17h"
